# Fix for "Analytics not working": the DAWNETTA_MBM_Worked log was missing
# its most recent ticket-assignment event (04-12-2023 17:27:42), which broke
# downstream reporting that expects every sheet to end on the latest event.
# Append the missing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DAWNETTA_MBM_Worked")

# Sheet currently ends at row 64 (A1:B64) -> new event goes to row 65.
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "17:27:42 04-12-2023"
$ws.Cells.Item($newRow, 2).Value = "Automatically Assigned Ticket"

# Widen column A a bit (as in the authored workbook) and reposition the
# view/selection near the bottom of the refreshed log.
$ws.Columns.Item(1).ColumnWidth = 32.45

$ws.Activate()
$ws.Range("A66:G68").Select() | Out-Null
